# "Add more analyzed papers"
#
# The Studies sheet tracks a literature-review coding pass: column G holds a
# numeric screening code and column H holds the matching reason label. Rows
# 150-154 were previously coded "New" (5 / "New" text); this pass reclassifies
# them (and codes a further batch of previously-unset rows, 155-195) with the
# updated "New Journal" label (and, where applicable, other reason labels).
# Two new reason labels are introduced: "New Journal" and "Unavailable".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Studies")

# Row -> (new G code [or $null if unchanged], new H reason label) taken from
# the reviewer's updated coding pass.
$updates = @(
    @{ Row = 150; G = $null; H = "New Journal" },
    @{ Row = 151; G = $null; H = "New Journal" },
    @{ Row = 152; G = $null; H = "New Journal" },
    @{ Row = 153; G = $null; H = "New Journal" },
    @{ Row = 154; G = $null; H = "New Journal" },
    @{ Row = 155; G = 5; H = "New Journal" },
    @{ Row = 156; G = 6; H = "Irrelevant focus" },
    @{ Row = 157; G = 5; H = "New Journal" },
    @{ Row = 158; G = 5; H = "New Journal" },
    @{ Row = 159; G = 6; H = "Irrelevant focus" },
    @{ Row = 160; G = 6; H = "Irrelevant focus" },
    @{ Row = 161; G = 6; H = "Unavailable" },
    @{ Row = 162; G = -1; H = "Duplicate" },
    @{ Row = 163; G = 6; H = "Theoretical paper" },
    @{ Row = 164; G = -1; H = "Duplicate" },
    @{ Row = 165; G = 6; H = "Irrelevant focus" },
    @{ Row = 166; G = 5; H = "New Journal" },
    @{ Row = 167; G = 5; H = "New Journal" },
    @{ Row = 168; G = 6; H = "Theoretical paper" },
    @{ Row = 169; G = 6; H = "Irrelevant focus" },
    @{ Row = 170; G = 6; H = "Theoretical paper" },
    @{ Row = 171; G = 6; H = "Theoretical paper" },
    @{ Row = 172; G = 6; H = "Irrelevant focus" },
    @{ Row = 173; G = 6; H = "Irrelevant focus" },
    @{ Row = 174; G = 6; H = "Theoretical paper" },
    @{ Row = 175; G = 6; H = "Theoretical paper" },
    @{ Row = 176; G = 6; H = "Irrelevant focus" },
    @{ Row = 177; G = 6; H = "Irrelevant focus" },
    @{ Row = 178; G = 6; H = "Theoretical paper" },
    @{ Row = 179; G = 6; H = "Theoretical paper" },
    @{ Row = 180; G = 6; H = "Theoretical paper" },
    @{ Row = 181; G = -1; H = "Duplicate" },
    @{ Row = 182; G = -1; H = "Duplicate" },
    @{ Row = 183; G = -1; H = "Duplicate" },
    @{ Row = 184; G = 5; H = "New Journal" },
    @{ Row = 185; G = -1; H = "Duplicate" },
    @{ Row = 186; G = -1; H = "Duplicate" },
    @{ Row = 187; G = -1; H = "Duplicate" },
    @{ Row = 188; G = -1; H = "Duplicate" },
    @{ Row = 189; G = 6; H = "Irrelevant focus" },
    @{ Row = 190; G = 6; H = "Theoretical paper" },
    @{ Row = 191; G = 6; H = "Irrelevant focus" },
    @{ Row = 192; G = -1; H = "Duplicate" },
    @{ Row = 193; G = 5; H = "New Journal" },
    @{ Row = 194; G = -1; H = "Duplicate" },
    @{ Row = 195; G = -1; H = "Duplicate" }
)

foreach ($u in $updates) {
    if ($null -ne $u.G) {
        $ws.Cells.Item($u.Row, 7).Value = $u.G
    }
    $ws.Cells.Item($u.Row, 8).Value = $u.H
}

# Reflect the reviewer's resulting scroll/selection position: they ended the
# pass around row 189, with the frozen pane scrolled to column F and the
# cursor on F196.
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 189
$aw.ScrollColumn = 6
$ws.Range("F196").Select()

$wb.Save()
